$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the 5th athlete's entry: "Co Ca" -> "Lê Quang Liêm", and fill in their unit as "HCM"
$ws.Range("B8").Value = "Lê Quang Liêm"
$ws.Range("C8").Value = "HCM"

# Remove the 7th athlete row entirely (STT, name, unit) - "Chim Sẻ Đi Nắng" / "Hà Nội"
$ws.Rows(10).Delete()

# Update the view: scroll back to the top and move the active selection
$ws.Activate()
$ws.Range("H6").Select()
